$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This weekly price sheet gets a new observation inserted at row 365; every
# existing row from 365 down to 406 shifts down by one (to 366..407).
$ws.Rows.Item(365).Insert()

# Seed the newly-inserted row 365 with the same fixed descriptive columns as
# its neighbours (market/region/category/quality/unit/origin/classification),
# then the new observation's own date, volume, prices and $/Kg.
$ws.Range("A365:R365").Value2 = $ws.Range("A366:R366").Value2

$ws.Range("D365").Value2 = 45194
$ws.Range("J365").Value2 = 80
$ws.Range("K365").Value2 = 4000
$ws.Range("L365").Value2 = 4000
$ws.Range("M365").Value2 = 4000
$ws.Range("P365").Value2 = 1333
